# Adds the "26-aug" column (BV) to the "Prix Spot" sheet, mirroring the
# automatic daily EPEX spot price export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Header cell, same bold/centered/bordered style as the other date headers.
$ws.Range("BU1").Copy()
$ws.Range("BV1").PasteSpecial(-4122)
$ws.Range("BV1").Value = "26-aug"

# Hourly price values for the new day.
$ws.Range("BV2").Value = 94.83
$ws.Range("BV3").Value = 86.78
$ws.Range("BV4").Value = 70.37
$ws.Range("BV5").Value = 63.58
$ws.Range("BV6").Value = 58.35
$ws.Range("BV7").Value = 66.86
$ws.Range("BV8").Value = 86.97
$ws.Range("BV9").Value = 95.44
$ws.Range("BV10").Value = 97.48
$ws.Range("BV11").Value = 90
$ws.Range("BV12").Value = 70
$ws.Range("BV13").Value = 60.65
$ws.Range("BV14").Value = 41.25
$ws.Range("BV15").Value = 26.28
$ws.Range("BV16").Value = 38.05
$ws.Range("BV17").Value = 47.44
$ws.Range("BV18").Value = 60
$ws.Range("BV19").Value = 65.03
$ws.Range("BV20").Value = 76.09999999999999
$ws.Range("BV21").Value = 103.14
$ws.Range("BV22").Value = 117
$ws.Range("BV23").Value = 123.46
$ws.Range("BV24").Value = 109.5
$ws.Range("BV25").Value = 97.90000000000001
